$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'" + '30.287.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'" + '2.091.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'" + '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.69%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("B5").Value = 'BNB'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'" + '341.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("B6").Value = 'USDC'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'" + '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").Value = 'XRP'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'" + '0.5309'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = 'Cardano'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'" + '0.4383'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").Value = 'OKB'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'" + '54.43'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'" + '0.09349'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = 'Polygon'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'" + '1.175'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = 'Solana'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'" + '24.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'" + '8.569'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.15%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'" + '6.887'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'" + '2.069.89'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.53%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'" + '101.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'" + '0.00001158'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'" + '1.003'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'" + '21.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = 'TRON'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'" + '0.06713'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'" + '6.353'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = 'Dai'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'" + '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'" + '30.260.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'" + '12.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'" + '2.308'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'" + '2.229.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.16%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'" + '21.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'" + '6.871'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.20%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = 'Monero'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'" + '162.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'" + '2.500'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'" + '133.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'" + '1.134'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'" + '1.669'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = 'Stellar'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'" + '0.1052'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'" + '6.276'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'" + '3.908'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = 'FraxShare'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'" + '10.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = 'VeChain'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'" + '0.02624'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = 'Hedera'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'" + '0.06761'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = 'Aptos'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'" + '12.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'" + '0.6955'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'" + '1.339'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = 'Algorand'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'" + '0.2212'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'" + '0.6767'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'" + '2.392'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'" + '14.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = 'Frax'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'" + '1.001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'" + '1.286'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.56%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'" + '3.630'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'" + '0.00000000350'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'" + '1.205'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.19%  '
$ws.Range("E51").Style = "Normal"
